# Generate Report for handoff
# Updates the status of the "2a4ef3df-2a97-4c83-90a1-9918c8ba7ad3.md" file
# to "Ready for handoff" on the Overview, zh-cn, and de-de sheets, and
# refreshes the "Latest Handoff Datetime" for the relevant rows.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D2").Value = "2016-01-26 06:29:55"
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-01-26 06:29:55"

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D2").Value = "2016-01-26 06:30:11"
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-01-26 06:30:11"
